$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay text; force text format so Excel
# does not reinterpret strings like "63.830.11" or "142.20" as numbers,
# then restore the default "Normal" style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.830.11"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.058.65"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "559.36"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "142.20"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.056.83"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "6.08"
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "35.29"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "3.565.81"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "63.917.42"
$ws.Range("D17").Value = "3.058.37"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "486.26"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "14.30"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "0.690"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "14.70"
$ws.Range("E23").Value = "  +8.44%  "
$ws.Range("D24").Value = "7.50"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "82.49"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "8.16"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "2.05"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "26.43"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "2.55"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").Value = "6.26"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").Value = "54.94"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "0.0411"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").Value = "441.91"
$ws.Range("E38").Value = "  -5.56%  "
$ws.Range("D39").Value = "0.0816"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "3.047.87"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "8.35"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  -9.12%  "
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("D44").Value = "0.275"
$ws.Range("E44").Value = "  +6.16%  "
$ws.Range("D45").Value = "27.90"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").Value = "2.24"
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "0.0₃0517"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").Value = "117.44"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "2.14"
$ws.Range("E51").Value = "  +2.98%  "

$priceRange.Style = "Normal"
